# Update NATMI TPM-based ligand-receptor edge statistics with newly
# computed values (re-run of script with new TPM input).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("G2").Value = 1.611874666666667
$ws.Range("H2").Value = 4.835624
$ws.Range("I2").Value = 0.06646895152072402
$ws.Range("J2").Value = 0.06646895152072402
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2936666666666667
$ws.Range("N2").Value = 0.881
$ws.Range("O2").Value = 0.009113820319201367
$ws.Range("P2").Value = 0.009113820319201367
$ws.Range("Q2").Value = 0.4733538604444445
$ws.Range("R2").Value = 4.260184744
$ws.Range("S2").Value = 0.0006057860809655851
$ws.Range("T2").Value = 0.0006057860809655851

# Row 3 (ECs -> FAPs)
$ws.Range("G3").Value = 1.611874666666667
$ws.Range("H3").Value = 4.835624
$ws.Range("I3").Value = 0.06646895152072402
$ws.Range("J3").Value = 0.06646895152072402
$ws.Range("O3").Value = 0.870405726797791
$ws.Range("P3").Value = 0.870405726797791
$ws.Range("Q3").Value = 45.20715753685067
$ws.Range("R3").Value = 406.864417831656
$ws.Range("S3").Value = 0.05785495605788293
$ws.Range("T3").Value = 0.05785495605788293

# Row 4 (ECs -> MuSCs)
$ws.Range("G4").Value = 1.611874666666667
$ws.Range("H4").Value = 4.835624
$ws.Range("I4").Value = 0.06646895152072402
$ws.Range("J4").Value = 0.06646895152072402
$ws.Range("O4").Value = 0.1204804528830076
$ws.Range("P4").Value = 0.1204804528830076
$ws.Range("Q4").Value = 6.257517208246223
$ws.Range("R4").Value = 56.317654874216
$ws.Range("S4").Value = 0.008008209381875504
$ws.Range("T4").Value = 0.008008209381875504

# Row 5 (FAPs -> ECs)
$ws.Range("I5").Value = 0.6469909869698216
$ws.Range("J5").Value = 0.6469909869698216
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.2936666666666667
$ws.Range("N5").Value = 0.881
$ws.Range("O5").Value = 0.009113820319201367
$ws.Range("P5").Value = 0.009113820319201367
$ws.Range("Q5").Value = 4.607499807777778
$ws.Range("R5").Value = 41.46749827
$ws.Range("S5").Value = 0.005896559603385707
$ws.Range("T5").Value = 0.005896559603385707

# Row 6 (FAPs -> FAPs)
$ws.Range("I6").Value = 0.6469909869698216
$ws.Range("J6").Value = 0.6469909869698216
$ws.Range("O6").Value = 0.870405726797791
$ws.Range("P6").Value = 0.870405726797791
$ws.Range("S6").Value = 0.5631446602450877
$ws.Range("T6").Value = 0.5631446602450877

# Row 7 (FAPs -> MuSCs)
$ws.Range("I7").Value = 0.6469909869698216
$ws.Range("J7").Value = 0.6469909869698216
$ws.Range("O7").Value = 0.1204804528830076
$ws.Range("P7").Value = 0.1204804528830076
$ws.Range("S7").Value = 0.07794976712134816
$ws.Range("T7").Value = 0.07794976712134816

# Row 8 (MuSCs -> ECs)
$ws.Range("I8").Value = 0.2865400615094543
$ws.Range("J8").Value = 0.2865400615094543
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.2936666666666667
$ws.Range("N8").Value = 0.881
$ws.Range("O8").Value = 0.009113820319201367
$ws.Range("P8").Value = 0.009113820319201367
$ws.Range("Q8").Value = 2.040574451444445
$ws.Range("R8").Value = 18.365170063
$ws.Range("S8").Value = 0.002611474634850074
$ws.Range("T8").Value = 0.002611474634850074

# Row 9 (MuSCs -> FAPs)
$ws.Range("I9").Value = 0.2865400615094543
$ws.Range("J9").Value = 0.2865400615094543
$ws.Range("O9").Value = 0.870405726797791
$ws.Range("P9").Value = 0.870405726797791
$ws.Range("S9").Value = 0.2494061104948203
$ws.Range("T9").Value = 0.2494061104948203

# Row 10 (MuSCs -> MuSCs)
$ws.Range("I10").Value = 0.2865400615094543
$ws.Range("J10").Value = 0.2865400615094543
$ws.Range("O10").Value = 0.1204804528830076
$ws.Range("P10").Value = 0.1204804528830076
$ws.Range("S10").Value = 0.0345224763797839
$ws.Range("T10").Value = 0.0345224763797839
